$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.989.20"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.47"
$ws.Range("E3").Value = "  -2.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.98"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4592"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3824"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07719"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9808"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.09"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.903.82"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.677"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.936"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07020"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.95"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009474"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.73"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.963.46"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.05"
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.668"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.74"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.851"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09272"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8658"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.072"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.251"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05750"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5519"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.434"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1757"
$ws.Range("E41").Value = "  -2.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.873"
$ws.Range("E42").Value = "  +4.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.331"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5182"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.23"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06845"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000002614"
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.063"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.99"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.784"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2871"
$ws.Range("E51").Value = "  -3.87%  "
